# prototype sql parsing (select command)
#
# Adds a new "what to do with the last token?" state column (A) to the
# state table, relabels column B from "index" to "state", and relabels
# column F from "value" to "string".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Range("A1").Value = "what to do with the last token?"
$ws.Range("B1").Value = "state"
$ws.Range("F1").Value = "string"

# --- Column A (new "what to do with the last token?" state values) --
# Row 2 (the "start" state) no longer carries a label in column A.
$ws.Range("A2").Value = $null

$ws.Range("A3").Value  = 'parse_tree["command"] = '
$ws.Range("A4").Value  = 'parse_tree["fields"] += '
$ws.Range("A5").Value  = 'parse_tree["fields"] += '
$ws.Range("A8").Value  = 'parse_tree["table_name"] ='
$ws.Range("A11").Value = 'parse_tree["where"] += '
$ws.Range("A12").Value = 'parse_tree["where"] += '
$ws.Range("A13").Value = 'parse_tree["where"] += '
$ws.Range("A14").Value = 'parse_tree["command"] = '
$ws.Range("A16").Value = 'parse_tree["table_name"] ='
$ws.Range("A18").Value = 'parse_tree["fields"] += '
$ws.Range("A20").Value = 'parse_tree["command"] = '
$ws.Range("A22").Value = 'parse_tree["table_name"] ='
